# Add a new "entities" row describing the new "new" / "New Metadata" table
# that was added to rd3_portal (commit: "added new portal table").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("entities")

$ws.Range("A4").Value = "rd3_portal_release"
$ws.Range("B4").Value = "new"
$ws.Range("C4").Value = "New Metadata"
$ws.Range("D4").Value = 'New metadata that was found in the "All patches" table'
$ws.Range("F4").Value = "rd3_portal_release_attrTmplate"
